$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FERNANDEZ VALDERAS ERNESTO ALI", 105),
    @("CONTRERAS VALDERRAMA JULIA ALEJANDRA", 89),
    @("VALLE MAGALLAN EDUAR", 86),
    @("HIDALGO CUBAS LUISA YVONE", 81),
    @("ROMERO CHANAME YOSSELY TRINIDAD", 81),
    @("ZAVALETA MANAY JORGE LUIS", 81),
    @("GUTIERREZ CARLOS TERESA DE JESUS", 77),
    @("CAMACHO LINARES JUDITH ARLETT", 76),
    @("BALLENA ESQUÉN ASTRID CAROLINA", 73),
    @("HUMPIRE CASTILLO IRWIN DEIMER", 73),
    @("SEVERINO AVALOS MARJORIE ISABEL", 73),
    @("SENADOR ARBOLEDA GIANCARLOS EXEBIO", 71),
    @("ZEVALLOS PACHECO ZOILA XIMENA", 70),
    @("SALAZAR FLORES ANA LIZETH", 10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
